$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-02 -> 2023-09-03, serial 45171 -> 45172) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 150 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 45172
    }
}
